$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Straight value updates (existing cells, same style retained automatically)
# ---------------------------------------------------------------------------
$ws.Range("T2").Value   = "maa://22742 (91.67), *maa://20791 (62.16)"
$ws.Range("T3").Value   = "maa://24617 (89.83), **maa://20790 (43.48), ***maa://37170 (16.67), maa://45854 (84.0)"
$ws.Range("D6").Value   = "maa://42407 (96.36)"
$ws.Range("A8").Value   = "更新日期：2025.03.07 13:21:00"
$ws.Range("X8").Value   = "maa://21411 (95.96)"
$ws.Range("P9").Value   = "maa://22736 (83.5)"
$ws.Range("X9").Value   = "maa://26223 (97.95)"
$ws.Range("D10").Value  = "***maa://25695 (18.62), ***maa://39951 (14.04), ***maa://34206 (19.23), ***maa://39243 (25.0), *maa://45271 (58.97)"
$ws.Range("T10").Value  = "maa://27395 (96.53), maa://22755 (87.83), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("T11").Value  = "maa://22747 (92.55), maa://22501 (97.78), maa://45521 (85.0)"
$ws.Range("X11").Value  = "maa://36713 (97.78)"
$ws.Range("AB11").Value = "maa://29912 (97.22), maa://22516 (88.37), *maa://20794 (52.24)"
$ws.Range("H12").Value  = "maa://21867 (90.0), **maa://45826 (33.33)"
$ws.Range("X12").Value  = "maa://22753 (91.06), *maa://21485 (75.35), maa://37962 (90.7)"
$ws.Range("P13").Value  = "maa://22676 (92.97), *maa://22583 (75.0), *maa://22500 (58.7)"
$ws.Range("AB14").Value = "maa://22764 (97.18)"
$ws.Range("D15").Value  = "*maa://22743 (78.04), maa://22734 (84.17), *maa://30808 (64.18), **maa://36048 (46.77), maa://45058 (93.33)"
$ws.Range("P15").Value  = "maa://24762 (90.36), *maa://22727 (70.0)"
$ws.Range("D20").Value  = "maa://21432 (90.29), maa://25198 (93.64), *maa://20795 (50.77), maa://36680 (91.18)"
$ws.Range("H22").Value  = "maa://25236 (95.74), **maa://21678 (48.94), **maa://22735 (42.86)"
$ws.Range("X22").Value  = "maa://21282 (98.61), *maa://37649 (65.52)"
$ws.Range("D23").Value  = "***maa://28036 (28.77), *maa://41753 (55.0)"
$ws.Range("L23").Value  = "maa://39756 (95.76), maa://39875 (94.44)"
$ws.Range("X23").Value  = "*maa://28503 (69.14)"
$ws.Range("X24").Value  = "maa://29988 (84.11), maa://23504 (93.22), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (76.92), ***maa://22815 (23.08)"
$ws.Range("D26").Value  = "maa://41802 (94.74)"
$ws.Range("L28").Value  = "maa://30770 (81.25)"
$ws.Range("T28").Value  = "*maa://29765 (64.71), maa://23263 (95.28)"
$ws.Range("D30").Value  = "maa://45792 (93.33)"
$ws.Range("L31").Value  = "maa://35926 (93.52), maa://36258 (84.75), *maa://43904 (72.73)"
$ws.Range("H32").Value  = "maa://21895 (97.04), maa://36667 (97.65), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("AF35").Value = "maa://39479 (89.47)"
$ws.Range("L37").Value  = "maa://45718 (97.96), *maa://47069 (73.33), maa://45789 (100.0)"
$ws.Range("P37").Value  = "maa://21280 (89.5), *maa://21239 (66.67)"
$ws.Range("P38").Value  = "*maa://24383 (68.93)"
$ws.Range("H39").Value  = "maa://36670 (89.11), maa://25199 (84.82), maa://30434 (91.67), ***maa://25036 (16.0), maa://45059 (83.33), *maa://44165 (66.67)"
$ws.Range("H58").Value  = "*maa://37964 (59.52)"
$ws.Range("H60").Value  = "*maa://40438 (69.84)"

# ---------------------------------------------------------------------------
# Row 34: O34/P34 go from "-" to "0" / "None".
# O34 keeps its General-formatted style, and "0" must stay a literal text
# value rather than collapse into the number 0 - write it as a text formula
# then freeze the result back into a plain value so no quote-prefix marker
# (which would change the cell style) is left behind.
# P34 is already Text-formatted (s="2"), so a normal value assignment keeps
# it a string automatically.
# ---------------------------------------------------------------------------
$ws.Range("O34").Formula = "=""0"""
$ws.Range("O34").Copy()
$ws.Range("O34").PasteSpecial(-4163)
$ws.Range("P34").Value = "None"

# ---------------------------------------------------------------------------
# New operator rows appended to existing groups. Each new "code" cell must
# take the same style as its row's other blank placeholder cells (style 1)
# even when the column default style is 2, so copy formats from a known
# style-1 neighbour before writing the text.
# ---------------------------------------------------------------------------

# Row 44: new 特种 (specialist) entry "钼铅"
$ws.Range("U44").Copy()
$ws.Range("AF44").PasteSpecial(-4122)
$ws.Range("U44").Copy()
$ws.Range("AG44").PasteSpecial(-4122)
$ws.Range("AD44").Value = "钼铅"
$ws.Range("AE44").Value = "-"
$ws.Range("AF44").Value = "-"

# Row 52: new 术师 (caster) entry "死芒"
$ws.Range("Q52").Copy()
$ws.Range("T52").PasteSpecial(-4122)
$ws.Range("Q52").Copy()
$ws.Range("U52").PasteSpecial(-4122)
$ws.Range("R52").Value = "死芒"
$ws.Range("S52").Value = "-"
$ws.Range("T52").Value = "-"

# Row 53: new 狙击 (sniper) entry "水灯心"
$ws.Range("I53").Copy()
$ws.Range("P53").PasteSpecial(-4122)
$ws.Range("I53").Copy()
$ws.Range("Q53").PasteSpecial(-4122)
$ws.Range("N53").Value = "水灯心"
$ws.Range("O53").Value = "-"
$ws.Range("P53").Value = "-"
